# Scheduled-runner market data refresh for the Brynhildr leve-profit tracker.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) on a handful of rows across the ALC/ARM/BSM/CRP/CUL/LTW/WVR
# sheets to reflect freshly-pulled market board prices. Some rows gain or
# lose a LeveProfit cell entirely depending on whether that side (NQ/HQ) of
# the leve has a value this run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 2385.6
$ws.Range("I10").Value = 3501.5
$ws.Range("J10").Value = 1641.6666
$ws.Range("K10").Value = 3501.5
$ws.Range("L10").Value = 1641.6666
$ws.Range("M10").Value = -3208.5
$ws.Range("N10").Value = -2227.6666

$ws.Range("H69").Value = 7114.9
$ws.Range("J69").Value = 7114.9
$ws.Range("L69").Value = 21344.7
$ws.Range("N69").Value = -23092.7

$ws.Range("H72").Value = 7114.9
$ws.Range("J72").Value = 7114.9
$ws.Range("L72").Value = 64034.1
$ws.Range("N72").Value = -72770.10000000001

$ws.Range("H106").Value = 7057.6665
$ws.Range("I106").Value = 7057.6665
$ws.Range("K106").Value = 7057.6665
$ws.Range("M106").Value = -6426.6665

$ws.Range("H138").Value = 3528.2942
$ws.Range("I138").Value = 2826
$ws.Range("K138").Value = 8478
$ws.Range("M138").Value = -3338

$ws.Range("H141").Value = 3663.2258
$ws.Range("I141").Value = 1751.2632
$ws.Range("K141").Value = 5253.7896
$ws.Range("M141").Value = -73.78960000000006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H31").Value = 10097.4
$ws.Range("I31").Value = 10097.4
$ws.Range("K31").Value = 10097.4
$ws.Range("M31").Value = -9803.4

$ws.Range("H32").Value = 129962.984
$ws.Range("I32").Value = 182768.39
$ws.Range("J32").Value = 20440.666
$ws.Range("K32").Value = 182768.39
$ws.Range("L32").Value = 20440.666
$ws.Range("M32").Value = -182481.39
$ws.Range("N32").Value = -21014.666

$ws.Range("H82").Value = 20181
$ws.Range("J82").Value = 20181
$ws.Range("L82").Value = 20181
$ws.Range("N82").Value = -20903

$ws.Range("H85").Value = 20181
$ws.Range("J85").Value = 20181
$ws.Range("L85").Value = 20181
$ws.Range("N85").Value = -22677

$ws.Range("H110").Value = 975
$ws.Range("I110").Value = 761.45
$ws.Range("K110").Value = 761.45
$ws.Range("M110").Value = 1283.55

$ws.Range("H132").Value = 718274.7
$ws.Range("I132").Value = 810448.9
$ws.Range("J132").Value = 3924.5
$ws.Range("K132").Value = 2431346.7
$ws.Range("L132").Value = 11773.5
$ws.Range("M132").Value = -2428816.7
$ws.Range("N132").Value = -16833.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H19").Value = 505

$ws.Range("H99").Value = 15135.091
$ws.Range("J99").Value = 1774.75
$ws.Range("L99").Value = 1774.75
$ws.Range("N99").Value = -4770.75

$ws.Range("H107").Value = 771.5
$ws.Range("I107").Value = 575.3043
$ws.Range("K107").Value = 575.3043
$ws.Range("M107").Value = 1344.6957

$ws.Range("H122").Value = 30000
$ws.Range("J122").Value = 30000
$ws.Range("L122").Value = 30000
$ws.Range("N122").Value = -39800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 51950
$ws.Range("J75").Value = 51950
$ws.Range("L75").Value = 51950
$ws.Range("N75").Value = -53946

$ws.Range("H78").Value = 51950
$ws.Range("J78").Value = 51950
$ws.Range("L78").Value = 155850
$ws.Range("N78").Value = -165834

$ws.Range("H99").Value = 24611.834
$ws.Range("I99").Value = 30671.285
$ws.Range("J99").Value = 16128.6
$ws.Range("K99").Value = 30671.285
$ws.Range("L99").Value = 16128.6
$ws.Range("M99").Value = -29173.285
$ws.Range("N99").Value = -19124.6

$ws.Range("H124").Value = 18625
$ws.Range("J124").Value = 18625
$ws.Range("L124").Value = 18625
$ws.Range("N124").Value = -23535

$ws.Range("H126").Value = 24611.834
$ws.Range("I126").Value = 30671.285
$ws.Range("J126").Value = 16128.6
$ws.Range("K126").Value = 92013.855
$ws.Range("L126").Value = 48385.8
$ws.Range("M126").Value = -89543.855
$ws.Range("N126").Value = -53325.8

$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 83333384
$ws.Range("I11").Value = 52.2
$ws.Range("J11").Value = 500000060
$ws.Range("K11").Value = 156.6
$ws.Range("L11").Value = 1500000180
$ws.Range("M11").Value = -16.60000000000002
$ws.Range("N11").Value = -1500000460

$ws.Range("H23").Value = 118.23077
$ws.Range("J23").Value = 138.625
$ws.Range("L23").Value = 415.875
$ws.Range("N23").Value = -885.875

$ws.Range("H42").Value = 994
$ws.Range("J42").Value = 994
$ws.Range("L42").Value = 2982
$ws.Range("N42").Value = -4050

$ws.Range("H62").Value = 2733.111
$ws.Range("J62").Value = 2666.6667
$ws.Range("L62").Value = 8000.000100000001
$ws.Range("N62").Value = -9372.000100000001

$ws.Range("H65").Value = 2733.111
$ws.Range("J65").Value = 2666.6667
$ws.Range("L65").Value = 24000.0003
$ws.Range("N65").Value = -30864.0003

$ws.Range("H68").Value = 6199.4814
$ws.Range("J68").Value = 8436.895
$ws.Range("L68").Value = 25310.685
$ws.Range("N68").Value = -26932.685

$ws.Range("H71").Value = 6199.4814
$ws.Range("J71").Value = 8436.895
$ws.Range("L71").Value = 75932.05500000001
$ws.Range("N71").Value = -84044.05500000001

$ws.Range("H82").Value = 15001.733
$ws.Range("I82").Value = 5013
$ws.Range("K82").Value = 15039
$ws.Range("M82").Value = -14633

$ws.Range("H85").Value = 15001.733
$ws.Range("I85").Value = 5013
$ws.Range("K85").Value = 15039
$ws.Range("M85").Value = -13635

$ws.Range("H92").Value = 271
$ws.Range("I92").Value = 246.33333
$ws.Range("J92").Value = 289.5
$ws.Range("K92").Value = 738.99999
$ws.Range("L92").Value = 868.5
$ws.Range("M92").Value = 509.00001
$ws.Range("N92").Value = -3364.5

$ws.Range("H97").Value = 595
$ws.Range("J97").Value = 595
$ws.Range("L97").Value = 1785
$ws.Range("N97").Value = -2777

$ws.Range("H128").Value = 191666.67
$ws.Range("I128").Value = 191666.67
$ws.Range("K128").Value = 575000.01
$ws.Range("M128").Value = -570020.01

$ws.Range("H137").Value = 8736.468999999999
$ws.Range("J137").Value = 10528.318
$ws.Range("L137").Value = 31584.954
$ws.Range("N137").Value = -41784.954

$ws.Range("H139").Value = 5701.846
$ws.Range("J139").Value = 11174.9
$ws.Range("L139").Value = 33524.7
$ws.Range("N139").Value = -43804.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 30253
$ws.Range("I99").Value = 30253
$ws.Range("K99").Value = 30253
$ws.Range("M99").Value = -27258

$ws.Range("H100").Value = 2322.4443
$ws.Range("I100").Value = 2322.4443
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2322.4443
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1781.4443
$ws.Range("N100").ClearContents()

$ws.Range("H133").Value = 70982.39999999999
$ws.Range("J133").Value = 70982.39999999999
$ws.Range("L133").Value = 70982.39999999999
$ws.Range("N133").Value = -76042.39999999999

$ws.Range("H136").Value = 14715756
$ws.Range("I136").Value = 15629192
$ws.Range("J136").Value = 13903813
$ws.Range("K136").Value = 46887576
$ws.Range("L136").Value = 41711439
$ws.Range("M136").Value = -46885026
$ws.Range("N136").Value = -41716539

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6946516
$ws.Range("I132").Value = 9261121
$ws.Range("J132").Value = 2700.6667
$ws.Range("K132").Value = 27783363
$ws.Range("L132").Value = 8102.000100000001
$ws.Range("M132").Value = -27780833
$ws.Range("N132").Value = -13162.0001

Write-Output "Brynhildr_Profits: updated market data on ALC/ARM/BSM/CRP/CUL/LTW/WVR"
